$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new timesheet entry for 1/22/2025
$ws.Range("A7").Value = [DateTime]"2025-01-22"
$ws.Range("C7").Value = "cleaning and sharing data, reviewing email text"
$ws.Range("B7").Value = ".5 hours"
$ws.Range("D7").Value = "Grant"

# Widen column C to accommodate the longer task text (matches Excel's
# auto best-fit sizing for the new, much longer entry in column C)
$ws.Columns.Item(3).ColumnWidth = 38.65

# Update selection to match post-edit state
$ws.Range("A8").Select()
